$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-AccountRow($account) {
    $found = $ws.Columns.Item(1).Find($account)
    if ($found -eq $null) {
        return -1
    }
    return $found.Row
}

function Remove-AccountRow($account) {
    $r = Find-AccountRow $account
    if ($r -ge 1) {
        $ws.Rows.Item($r).Delete()
    }
}

function Insert-RowBeforeAccount($beforeAccount, $account, $name, $saldo) {
    $r = Find-AccountRow $beforeAccount
    if ($r -lt 1) {
        throw "Anchor account not found: $beforeAccount"
    }
    $ws.Rows.Item($r).Insert()
    # Force text format on the account-number cell so leading zeros survive
    # (Conta values are zero-padded account numbers, stored as text).
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $account
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $saldo
}

# --- Remove rows that no longer exist / whose old value will be replaced ---
# ERIK (005995120) account was closed / removed entirely.
Remove-AccountRow "005995120"
# LUCIANO (005002390) account was closed / removed entirely.
Remove-AccountRow "005002390"
# AHMAD (004368468) old balance row removed; re-added below with new balance.
Remove-AccountRow "004368468"
# GUSTAVO (005591536) old balance row removed; re-added below with new balance.
Remove-AccountRow "005591536"
# TIAGO (005924958) old balance row removed; re-added below with new balance.
Remove-AccountRow "005924958"

# --- Insert rows with new/updated balances, keeping the descending Saldo sort ---
# NATALIA is a new account with balance 5000, placed right after LARISSA (5000).
Insert-RowBeforeAccount "004398253" "004482102" "NATALIA" 5000
# AHMAD's updated balance (1721.23) now sits just after NATALIA, before EULER.
Insert-RowBeforeAccount "004398253" "004368468" "AHMAD" 1721.23

# NILSON is a new account with balance 492.2, placed before CLAUDIA (461.89).
Insert-RowBeforeAccount "005044389" "008115273" "NILSON" 492.2

# GUSTAVO's updated balance (447.36) now sits before ANA (446.18).
Insert-RowBeforeAccount "004432579" "005591536" "GUSTAVO" 447.36

# TIAGO's updated balance (438.4) now sits before MARCO (365.23).
Insert-RowBeforeAccount "004436055" "005924958" "TIAGO" 438.4
